$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.863.88"
$ws.Range("D3").Value = "1.563.59"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.77"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.79"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0864"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "1.787.30"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "1.561.44"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "26.884.38"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.31"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.34"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.37"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").Value = "0.0₃0682"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.21"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.68"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0466"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("E31").Value = "  -3.37%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "1.392.73"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("E35").Value = "  -0.95%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.924"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.17%  "
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.531"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.66%  "
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.992"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.52"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.79"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.88"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").Value = "1.700.73"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.73"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0503"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0973"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0953"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.48%  "
